$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append three new daily rows (27/28/29-Nov-2020) to the COVID patient
# condition tracking table.
$ws.Range("A258").Value = 44162
$ws.Range("B258").Value = 1341
$ws.Range("C258").Value = 453
$ws.Range("D258").Value = 433
$ws.Range("E258").Value = 128
$ws.Range("F258").Value = 14

$ws.Range("A259").Value = 44163
$ws.Range("B259").Value = 1405
$ws.Range("C259").Value = 379
$ws.Range("D259").Value = 413
$ws.Range("E259").Value = 138
$ws.Range("F259").Value = 11

$ws.Range("A260").Value = 44164
$ws.Range("B260").Value = 1313
$ws.Range("C260").Value = 375
$ws.Range("D260").Value = 434
$ws.Range("E260").Value = 128
$ws.Range("F260").Value = 11

# Match the formatting (date style in column A, centered numbers in B:F)
# used by the rest of the table by copying it down from the last existing
# row.
$ws.Range("A257:F257").Copy()
$ws.Range("A258:F260").PasteSpecial(-4122)  # xlPasteFormats

# Grow the Excel Table/AutoFilter so the new rows are included.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F260"))

# Move the active selection to just past the new last row, as in the
# source workbook.
$ws.Range("F261").Select() | Out-Null
